# Update on 2018-1-2, 支出生活费300
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet ("2017年" -> "Summary"); this also updates the
#    _xlnm._FilterDatabase defined name which refers to the sheet by name.
$ws.Name = "Summary"

# 2) Fill in row 27 (record #25): a new "生活费" (living expenses) expense
#    of 300 on 2018-1-2 (serial 43102), matching the style/formatting of
#    the row above it (row 26).
$ws.Range("B26:G26").Copy()
$ws.Range("B27:G27").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C27").Value = "支出"
$ws.Range("D27").Value = 300
$ws.Range("E27").Value = 43102
$ws.Range("F27").Value = "生活费"
$ws.Range("G27").Value = "生活费(1/1-1/10)"

# 3) Move the active selection to reflect where the user ended up editing.
$ws.Range("G32").Select()
